# Update calibration data with new costs (Brazil SE calibrated inputs)
# Rows (1-indexed as in the sheet):
#   Row 2  -> gdp_mmm_usd                              (J2:AS2)
#   Row 8  -> elasticity_gnrl_rate_occupancy_to_gdppc   (J8:AS8)
#   Row 9  -> frac_gnrl_eating_red_meat                 (J9:AS9)
#   Row 13 -> occrateinit_gnrl_occupancy                (J13:AS13)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: J2:AS2 ---
$row2 = New-Object 'object[,]' 1,36
$row2[0,0]  = 2431.942902
$row2[0,1]  = 2421.7722246
$row2[0,2]  = 2411.6015472
$row2[0,3]  = 2401.4308698
$row2[0,4]  = 2391.2601924
$row2[0,5]  = 2381.089515
$row2[0,6]  = 2420.084159
$row2[0,7]  = 2459.078803
$row2[0,8]  = 2498.073447
$row2[0,9]  = 2537.068091
$row2[0,10] = 2576.062735
$row2[0,11] = 2644.5606574
$row2[0,12] = 2713.0585798
$row2[0,13] = 2781.5565022
$row2[0,14] = 2850.0544246
$row2[0,15] = 2918.552347
$row2[0,16] = 2984.2074356
$row2[0,17] = 3049.8625242
$row2[0,18] = 3115.5176128
$row2[0,19] = 3181.1727014
$row2[0,20] = 3246.82779
$row2[0,21] = 3314.4815632
$row2[0,22] = 3382.1353364
$row2[0,23] = 3449.7891096
$row2[0,24] = 3517.4428828
$row2[0,25] = 3585.096656
$row2[0,26] = 3653.2644508
$row2[0,27] = 3721.4322456
$row2[0,28] = 3789.6000404
$row2[0,29] = 3857.7678352
$row2[0,30] = 3925.93563
$row2[0,31] = 3993.1127378
$row2[0,32] = 4060.2898456
$row2[0,33] = 4127.4669534
$row2[0,34] = 4194.6440612
$row2[0,35] = 4261.821169
$ws.Range("J2:AS2").Value = $row2

# --- Row 8: J8:AS8 (constant -0.1 across all columns) ---
$row8 = New-Object 'object[,]' 1,36
for ($i = 0; $i -lt 36; $i++) { $row8[0,$i] = -0.1 }
$ws.Range("J8:AS8").Value = $row8

# --- Row 9: J9:AS9 (constant 1 across all columns) ---
$row9 = New-Object 'object[,]' 1,36
for ($i = 0; $i -lt 36; $i++) { $row9[0,$i] = 1 }
$ws.Range("J9:AS9").Value = $row9

# --- Row 13: J13:AS13 (constant 3.145207224 across all columns) ---
$row13 = New-Object 'object[,]' 1,36
for ($i = 0; $i -lt 36; $i++) { $row13[0,$i] = 3.145207224 }
$ws.Range("J13:AS13").Value = $row13
